# Add donor_organism.disease_profile COPD fields to the "Donor organism" sheet
# (template/hca_lung_template.xlsx) right after the "KCO PERCENT OF PREDICTED"
# column, and flip the active tab from "Donor organism" back to "Project".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Donor organism")

# Insert six new blank columns right before the old CJ (after CI, which holds
# "KCO PERCENT OF PREDICTED"). This shifts every column from CJ onward six
# places to the right (CJ->CP ... CY->DE), and inherits formatting from the
# column immediately to the left, matching the target layout.
$ws.Columns("CJ:CO").Insert()

# --- Row 1: field labels -------------------------------------------------
$ws.Range("CJ1").Value = "COPD - GOLD STAGE"
$ws.Range("CK1").Value = "COPD - MMRC GRADE"
$ws.Range("CL1").Value = "COPD - CAT SCORE"
$ws.Range("CM1").Value = "COPD - GOLD ABE ASSESSMENT"
$ws.Range("CN1").Value = "COPD PHENOTYPE"
$ws.Range("CO1").Value = "COPD - PERCENTAGE OF EMPHYSEMA"

# --- Row 2: field descriptions --------------------------------------------
$ws.Range("CJ2").Value = "Indicate the current GOLD stage (Global Initiative for Chronic Obstructive Lung Disease)."
$ws.Range("CK2").Value = "Indicate the Modified British Medical Research Council (mMRC) dyspnea scale grade"
$ws.Range("CL2").Value = "Indicate the COPD Assessment Test (CAT) score."
$ws.Range("CM2").Value = "Indicate the Global Initiative for Chronic Obstructive Lung Disease (GOLD) A, B, C, D assessment group if available."
$ws.Range("CN2").Value = "Indicate the COPD disease phenotype(s) of donor. Please indicate all applicable phenotypes of donor."
$ws.Range("CO2").Value = "Indicate the percentage of the lung that is affected by emphysema as judged based on non-invasive imaging, such as from a CT scan."

# --- Row 3: guidance / examples -------------------------------------------
$ws.Range("CJ3").Value = "Should be one of 1, 2, 3, 4"
$ws.Range("CK3").Value = "Should be one of 0, 1, 2, 3, 4"
$ws.Range("CL3").Value = "Should be between 0 and 40"
$ws.Range("CM3").Value = "Should be one of A, B, E"
$ws.Range("CN3").Value = "Should be one or more of: COPD not otherwise specified, COPD with emphysema, COPD with bronchitis, COPD with history of asthma, COPD with >300 eos in blood, COPD with allergy, COPD with Chronic Mucus Hypersecretion, COPD with frequent exacerbations"
$ws.Range("CO3").Value = "For example: 93; 85; 77"

# --- Row 4: programmatic field names --------------------------------------
$ws.Range("CJ4").Value = "donor_organism.disease_profile.copd_gold_stage"
$ws.Range("CK4").Value = "donor_organism.disease_profile.copd_mmrc_grade"
$ws.Range("CL4").Value = "donor_organism.disease_profile.copd_cat_score"
$ws.Range("CM4").Value = "donor_organism.disease_profile.copd_gold_abe_assessment"
$ws.Range("CN4").Value = "donor_organism.disease_profile.copd_phenotype"
$ws.Range("CO4").Value = "donor_organism.disease_profile.copd_emphysema_percentage"

# --- Column widths for the newly inserted columns -------------------------
# (target raw xlsx widths are 16.6640625 for CJ:CN and 19.5 for CO; the
# closest values reachable through the ColumnWidth property's own pixel
# rounding are used here)
$ws.Range("CJ1:CN1").ColumnWidth = 15.83
$ws.Range("CO1").ColumnWidth = 18.65

# --- Switch the active sheet back to "Project" ----------------------------
$project = $wb.Worksheets.Item("Project")
$project.Activate()
